{"js": "// Rewrite the memory-bound paragraph and the arithmetic-intensity paragraph\n// to match the author's edit (drop the Numpy mention, switch the matrix size\n// from 2000x2000 to 400x400, rephrase the multiply-cost sentence, and extend\n// the processor-utilization sentence at the end of the following paragraph).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet memParagraph = null;\nlet archParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t.indexOf(\"We are using random arrays generated by Numpy\") !== -1) {\n    memParagraph = paragraphs.items[i];\n  }\n  if (t.indexOf(\"For the high arithmetic intensity\") !== -1) {\n    archParagraph = paragraphs.items[i];\n  }\n}\n\nif (!memParagraph) {\n  throw new Error(\"Could not find the memory-bound paragraph\");\n}\nif (!archParagraph) {\n  throw new Error(\"Could not find the arithmetic-intensity paragraph\");\n}\n\n// --- Paragraph 1: memory-bound benchmark description -----------------------\n\n// Drop the \"We are using random arrays generated by Numpy.  \" lead-in.\nlet hits = memParagraph.search(\"We are using random arrays generated by Numpy.  For the mem\", { matchCase: true });\nhits.load(\"items/text\");\nawait context.sync();\nif (hits.items.length === 0) {\n  throw new Error(\"Could not find the Numpy lead-in sentence\");\n}\nhits.items[0].insertText(\"For the mem\", \"Replace\");\nawait context.sync();\n\n// 4 2000x2000 matrices -> 4 400x400 matrices\nhits = memParagraph.search(\"4 2000x2000 matrices\", { matchCase: true });\nhits.load(\"items/text\");\nawait context.sync();\nif (hits.items.length === 0) {\n  throw new Error(\"Could not find the matrix-size sentence\");\n}\nhits.items[0].insertText(\"4 400x400 matrices\", \"Replace\");\nawait context.sync();\n\n// Rewrite the \"multiplied pairwise ... summation of 2000 products\" sentence.\nconst oldMultiplySentence =\n  \"They are then multiplied pairwise, and the result discarded.   This is repeated until time has passed.  \" +\n  \"Each matrix multiply requires 2000 multiplications per row and then a summation of 2000 products.  \" +\n  \"Therefore, each matrix m\";\nconst newMultiplySentence =\n  \"They are then na\u00efvely matrix multiplied, such that each matrix multiply requires 400 multiplications per row and then a \" +\n  \"summation of 400 products for each of 400.  Therefore, each matrix m\";\nhits = memParagraph.search(oldMultiplySentence, { matchCase: true });\nhits.load(\"items/text\");\nawait context.sync();\nif (hits.items.length === 0) {\n  throw new Error(\"Could not find the multiply-cost sentence\");\n}\nhits.items[0].insertText(newMultiplySentence, \"Replace\");\nawait context.sync();\n\n// --- Paragraph 2: arithmetic-intensity / processor-utilization sentence ----\n\nconst oldUtilSentence =\n  \"100% processor utilization, we also employ 4 iterations of the python program itself.  \" +\n  \"This pushes the processor load to \";\nconst newUtilSentence =\n  \"100% processor utilization, but instead is about 35-40%, we also employ 3 iterations of the python program itself.  \" +\n  \"This pushes the processor load to 90-100% for the full 10 minutes, at least on the laptop.  \";\nhits = archParagraph.search(oldUtilSentence, { matchCase: true });\nhits.load(\"items/text\");\nawait context.sync();\nif (hits.items.length === 0) {\n  throw new Error(\"Could not find the processor-utilization sentence\");\n}\nhits.items[0].insertText(newUtilSentence, \"Replace\");\nawait context.sync();\n", "ps1": "# Rewrite the memory-bound paragraph and the arithmetic-intensity paragraph\n# to match the author's edit (drop the Numpy mention, switch the matrix size\n# from 2000x2000 to 400x400, rephrase the multiply-cost sentence, and extend\n# the processor-utilization sentence at the end of the following paragraph).\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Text = $findText\n    $rng.Find.Replacement.Text = $replaceText\n    $rng.Find.Forward = $true\n    $rng.Find.Wrap = 1\n    $rng.Find.MatchCase = $true\n    $rng.Find.MatchWholeWord = $false\n    $rng.Find.MatchWildcards = $false\n    # wdReplaceAll = 2\n    $rng.Find.Execute([ref]$findText, [ref]$true, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$true, [ref]$null, [ref]$null, [ref]$replaceText, [ref]2)\n}\n\n# Drop the \"We are using random arrays generated by Numpy.  \" lead-in.\nReplace-Text \"We are using random arrays generated by Numpy.  For the mem\" \"For the mem\"\n\n# 4 2000x2000 matrices -> 4 400x400 matrices\nReplace-Text \"4 2000x2000 matrices\" \"4 400x400 matrices\"\n\n# Rewrite the \"multiplied pairwise ... summation of 2000 products\" sentence.\nReplace-Text \"They are then multiplied pairwise, and the result discarded.   This is repeated until time has passed.  Each matrix multiply requires 2000 multiplications per row and then a summation of 2000 products.  Therefore, each matrix m\" \"They are then na\u00efvely matrix multiplied, such that each matrix multiply requires 400 multiplications per row and then a summation of 400 products for each of 400.  Therefore, each matrix m\"\n\n# Extend the processor-utilization sentence.\nReplace-Text \"100% processor utilization, we also employ 4 iterations of the python program itself.  This pushes the processor load to \" \"100% processor utilization, but instead is about 35-40%, we also employ 3 iterations of the python program itself.  This pushes the processor load to 90-100% for the full 10 minutes, at least on the laptop.  \"\n"}
